# Edit script for dheeraj_chand resume docx
# 1. Collapse the three CORE COMPETENCIES paragraphs into a single condensed
#    paragraph (category titles only, joined by bullets).
# 2. Append a new "TECHNICAL SKILLS" section (Heading2) at the end of the
#    document containing the previously-removed detail, reformatted with
#    semicolon separators and without the sub-bulleted tool lists.

$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------
# Step 1: Replace the CORE COMPETENCIES detail paragraphs
# ---------------------------------------------------------------------
# Find the first of the three long paragraphs and replace its text with
# the short summary line, then remove the following two paragraphs
# entirely (text + paragraph mark).

$found = $d.Content.Find.Execute(
    "Survey Methodology & Research Design: Survey Design*Expert Testimony and Consultation on Research Methodology",
    $false, $false, $true, $false, $false, $true, 1, $false,
    "Survey Methodology & Research Design " + $bullet + " Redistricting & Geospatial Analysis " + $bullet + " Data Analysis & Visualization",
    2
)

# Locate the paragraph that now holds the condensed Core Competencies line
# so we can remove the two now-redundant detail paragraphs that follow it.
$coreCompIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Survey Methodology & Research Design " + $bullet)) {
        $coreCompIndex = $i
        break
    }
}

$nextPara = $d.Paragraphs.Item($coreCompIndex + 1)
$afterPara = $d.Paragraphs.Item($coreCompIndex + 2)
$removeRange = $d.Range($nextPara.Range.Start, $afterPara.Range.End)
$removeRange.Delete()

# ---------------------------------------------------------------------
# Step 2: Append the new TECHNICAL SKILLS section at the end of the doc
# ---------------------------------------------------------------------
# Insert all four new (still-empty) paragraphs first, while the
# insertion point's contextual style is still "Normal" (inherited from
# the last "Built comprehensive survey..." bullet). Only afterwards do
# we stamp the Heading2 style onto the section title paragraph, and fill
# in the run text for every paragraph. This avoids the Heading2 style
# leaking into the three detail paragraphs that follow the heading.

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($n - 3)
$p1 = $d.Paragraphs.Item($n - 2)
$p2 = $d.Paragraphs.Item($n - 1)
$p3 = $d.Paragraphs.Item($n)

$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = "Heading2"

$p1.Range.Text = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Focus Groups and Qualitative Research Methodologies; Meta-analytical Dataset Development for Longitudinal Analysis; Survey Instrument Standardization and Call Methods Optimization; Expert Testimony and Consultation on Research Methodology"

$p2.Range.Text = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Spatial Clustering and Boundary Estimation without ML Requirements; Census Data Integration and Demographic Mapping; Court Case Analysis and Expert Testimony for Redistricting; Multi-tenant Data Warehouse Design for Electoral Analytics"

$p3.Range.Text = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Machine Learning and Predictive Modeling for Targeting; Big Data Analytics; Fraud Detection and Entity Resolution Systems; Multi-million Dollar Research Project Management"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
